$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('A2').Value2 = 'ECs'
$ws.Range('B2').Value2 = 'Ncam1'
$ws.Range('C2').Value2 = 'Ptprz1'
$ws.Range('D2').Value2 = 'ECs'
$ws.Range('E2').Value2 = [double]"3"
$ws.Range('F2').Value2 = [double]"1"
$ws.Range('G2').Value2 = [double]"0.9949870000000001"
$ws.Range('H2').Value2 = [double]"2.984961"
$ws.Range('I2').Value2 = [double]"0.03855738270564991"
$ws.Range('J2').Value2 = [double]"0.03855738270564991"
$ws.Range('K2').Value2 = [double]"1"
$ws.Range('L2').Value2 = [double]"0.3333333333333333"
$ws.Range('M2').Value2 = [double]"0.04936366666666667"
$ws.Range('N2').Value2 = [double]"0.148091"
$ws.Range('O2').Value2 = [double]"0.04616170608573571"
$ws.Range('P2').Value2 = [double]"0.0461617060857357"
$ws.Range('Q2').Value2 = [double]"0.04911620660566667"
$ws.Range('R2').Value2 = [double]"0.4420458594510001"
$ws.Range('S2').Value2 = [double]"0.00177987456789344"
$ws.Range('T2').Value2 = [double]"0.00177987456789344"

# Row 3
$ws.Range('A3').Value2 = 'ECs'
$ws.Range('B3').Value2 = 'Ncam1'
$ws.Range('C3').Value2 = 'Ptprz1'
$ws.Range('D3').Value2 = 'FAPs'
$ws.Range('E3').Value2 = [double]"3"
$ws.Range('F3').Value2 = [double]"1"
$ws.Range('G3').Value2 = [double]"0.9949870000000001"
$ws.Range('H3').Value2 = [double]"2.984961"
$ws.Range('I3').Value2 = [double]"0.03855738270564991"
$ws.Range('J3').Value2 = [double]"0.03855738270564991"
$ws.Range('K3').Value2 = [double]"1"
$ws.Range('L3').Value2 = [double]"0.3333333333333333"
$ws.Range('M3').Value2 = [double]"0.008616666666666667"
$ws.Range('N3').Value2 = [double]"0.02585"
$ws.Range('O3').Value2 = [double]"0.008057748967298944"
$ws.Range('P3').Value2 = [double]"0.008057748967298944"
$ws.Range('Q3').Value2 = [double]"0.008573471316666668"
$ws.Range('R3').Value2 = [double]"0.07716124185000001"
$ws.Range('S3').Value2 = [double]"0.0003106857106782007"
$ws.Range('T3').Value2 = [double]"0.0003106857106782007"

# Row 4
$ws.Range('A4').Value2 = 'ECs'
$ws.Range('B4').Value2 = 'Ncam1'
$ws.Range('C4').Value2 = 'Ptprz1'
$ws.Range('D4').Value2 = 'MuSCs'
$ws.Range('E4').Value2 = [double]"3"
$ws.Range('F4').Value2 = [double]"1"
$ws.Range('G4').Value2 = [double]"0.9949870000000001"
$ws.Range('H4').Value2 = [double]"2.984961"
$ws.Range('I4').Value2 = [double]"0.03855738270564991"
$ws.Range('J4').Value2 = [double]"0.03855738270564991"
$ws.Range('K4').Value2 = [double]"3"
$ws.Range('L4').Value2 = [double]"1"
$ws.Range('M4').Value2 = [double]"1.011383666666666"
$ws.Range('N4').Value2 = [double]"3.034151"
$ws.Range('O4').Value2 = [double]"0.9457805449469654"
$ws.Range('P4').Value2 = [double]"0.9457805449469653"
$ws.Range('Q4').Value2 = [double]"1.006313600345667"
$ws.Range('R4').Value2 = [double]"9.056822403110999"
$ws.Range('S4').Value2 = [double]"0.03646682242707827"
$ws.Range('T4').Value2 = [double]"0.03646682242707827"

# Row 5
$ws.Range('A5').Value2 = 'FAPs'
$ws.Range('B5').Value2 = 'Ncam1'
$ws.Range('C5').Value2 = 'Ptprz1'
$ws.Range('D5').Value2 = 'ECs'
$ws.Range('E5').Value2 = [double]"3"
$ws.Range('F5').Value2 = [double]"1"
$ws.Range('G5').Value2 = [double]"8.469728666666667"
$ws.Range('H5').Value2 = [double]"25.409186"
$ws.Range('I5').Value2 = [double]"0.3282159160005915"
$ws.Range('J5').Value2 = [double]"0.3282159160005916"
$ws.Range('K5').Value2 = [double]"1"
$ws.Range('L5').Value2 = [double]"0.3333333333333333"
$ws.Range('M5').Value2 = [double]"0.04936366666666667"
$ws.Range('N5').Value2 = [double]"0.148091"
$ws.Range('O5').Value2 = [double]"0.04616170608573571"
$ws.Range('P5').Value2 = [double]"0.0461617060857357"
$ws.Range('Q5').Value2 = [double]"0.4180968626584444"
$ws.Range('R5').Value2 = [double]"3.762871763926"
$ws.Range('S5').Value2 = [double]"0.01515100664707983"
$ws.Range('T5').Value2 = [double]"0.01515100664707983"

# Row 6
$ws.Range('A6').Value2 = 'FAPs'
$ws.Range('B6').Value2 = 'Ncam1'
$ws.Range('C6').Value2 = 'Ptprz1'
$ws.Range('D6').Value2 = 'FAPs'
$ws.Range('E6').Value2 = [double]"3"
$ws.Range('F6').Value2 = [double]"1"
$ws.Range('G6').Value2 = [double]"8.469728666666667"
$ws.Range('H6').Value2 = [double]"25.409186"
$ws.Range('I6').Value2 = [double]"0.3282159160005915"
$ws.Range('J6').Value2 = [double]"0.3282159160005916"
$ws.Range('K6').Value2 = [double]"1"
$ws.Range('L6').Value2 = [double]"0.3333333333333333"
$ws.Range('M6').Value2 = [double]"0.008616666666666667"
$ws.Range('N6').Value2 = [double]"0.02585"
$ws.Range('O6').Value2 = [double]"0.008057748967298944"
$ws.Range('P6').Value2 = [double]"0.008057748967298944"
$ws.Range('Q6').Value2 = [double]"0.07298082867777778"
$ws.Range('R6').Value2 = [double]"0.6568274581000001"
$ws.Range('S6').Value2 = [double]"0.002644681458204843"
$ws.Range('T6').Value2 = [double]"0.002644681458204844"

# Row 7
$ws.Range('A7').Value2 = 'FAPs'
$ws.Range('B7').Value2 = 'Ncam1'
$ws.Range('C7').Value2 = 'Ptprz1'
$ws.Range('D7').Value2 = 'MuSCs'
$ws.Range('E7').Value2 = [double]"3"
$ws.Range('F7').Value2 = [double]"1"
$ws.Range('G7').Value2 = [double]"8.469728666666667"
$ws.Range('H7').Value2 = [double]"25.409186"
$ws.Range('I7').Value2 = [double]"0.3282159160005915"
$ws.Range('J7').Value2 = [double]"0.3282159160005916"
$ws.Range('K7').Value2 = [double]"3"
$ws.Range('L7').Value2 = [double]"1"
$ws.Range('M7').Value2 = [double]"1.011383666666666"
$ws.Range('N7').Value2 = [double]"3.034151"
$ws.Range('O7').Value2 = [double]"0.9457805449469654"
$ws.Range('P7').Value2 = [double]"0.9457805449469653"
$ws.Range('Q7').Value2 = [double]"8.566145234565109"
$ws.Range('R7').Value2 = [double]"77.095307111086"
$ws.Range('S7').Value2 = [double]"0.3104202278953069"
$ws.Range('T7').Value2 = [double]"0.3104202278953069"

# Row 8
$ws.Range('A8').Value2 = 'Inflammatory-Mac'
$ws.Range('B8').Value2 = 'Ncam1'
$ws.Range('C8').Value2 = 'Ptprz1'
$ws.Range('D8').Value2 = 'ECs'
$ws.Range('E8').Value2 = [double]"2"
$ws.Range('F8').Value2 = [double]"0.6666666666666666"
$ws.Range('G8').Value2 = [double]"0.04495399999999999"
$ws.Range('H8').Value2 = [double]"0.134862"
$ws.Range('I8').Value2 = [double]"0.001742041435867791"
$ws.Range('J8').Value2 = [double]"0.001742041435867791"
$ws.Range('K8').Value2 = [double]"1"
$ws.Range('L8').Value2 = [double]"0.3333333333333333"
$ws.Range('M8').Value2 = [double]"0.04936366666666667"
$ws.Range('N8').Value2 = [double]"0.148091"
$ws.Range('O8').Value2 = [double]"0.04616170608573571"
$ws.Range('P8').Value2 = [double]"0.0461617060857357"
$ws.Range('Q8').Value2 = [double]"0.002219094271333333"
$ws.Range('R8').Value2 = [double]"0.019971848442"
$ws.Range('S8').Value2 = [double]"8.041560475170197E-05"
$ws.Range('T8').Value2 = [double]"8.041560475170197E-05"

# Row 9
$ws.Range('A9').Value2 = 'Inflammatory-Mac'
$ws.Range('B9').Value2 = 'Ncam1'
$ws.Range('C9').Value2 = 'Ptprz1'
$ws.Range('D9').Value2 = 'FAPs'
$ws.Range('E9').Value2 = [double]"2"
$ws.Range('F9').Value2 = [double]"0.6666666666666666"
$ws.Range('G9').Value2 = [double]"0.04495399999999999"
$ws.Range('H9').Value2 = [double]"0.134862"
$ws.Range('I9').Value2 = [double]"0.001742041435867791"
$ws.Range('J9').Value2 = [double]"0.001742041435867791"
$ws.Range('K9').Value2 = [double]"1"
$ws.Range('L9').Value2 = [double]"0.3333333333333333"
$ws.Range('M9').Value2 = [double]"0.008616666666666667"
$ws.Range('N9').Value2 = [double]"0.02585"
$ws.Range('O9').Value2 = [double]"0.008057748967298944"
$ws.Range('P9').Value2 = [double]"0.008057748967298944"
$ws.Range('Q9').Value2 = [double]"0.0003873536333333333"
$ws.Range('R9').Value2 = [double]"0.0034861827"
$ws.Range('S9').Value2 = [double]"1.403693258085566E-05"
$ws.Range('T9').Value2 = [double]"1.403693258085566E-05"

# Row 10
$ws.Range('A10').Value2 = 'Inflammatory-Mac'
$ws.Range('B10').Value2 = 'Ncam1'
$ws.Range('C10').Value2 = 'Ptprz1'
$ws.Range('D10').Value2 = 'MuSCs'
$ws.Range('E10').Value2 = [double]"2"
$ws.Range('F10').Value2 = [double]"0.6666666666666666"
$ws.Range('G10').Value2 = [double]"0.04495399999999999"
$ws.Range('H10').Value2 = [double]"0.134862"
$ws.Range('I10').Value2 = [double]"0.001742041435867791"
$ws.Range('J10').Value2 = [double]"0.001742041435867791"
$ws.Range('K10').Value2 = [double]"3"
$ws.Range('L10').Value2 = [double]"1"
$ws.Range('M10').Value2 = [double]"1.011383666666666"
$ws.Range('N10').Value2 = [double]"3.034151"
$ws.Range('O10').Value2 = [double]"0.9457805449469654"
$ws.Range('P10').Value2 = [double]"0.9457805449469653"
$ws.Range('Q10').Value2 = [double]"0.04546574135133332"
$ws.Range('R10').Value2 = [double]"0.4091916721619999"
$ws.Range('S10').Value2 = [double]"0.001647588898535233"
$ws.Range('T10').Value2 = [double]"0.001647588898535233"

# Row 11
$ws.Range('A11').Value2 = 'MuSCs'
$ws.Range('B11').Value2 = 'Ncam1'
$ws.Range('C11').Value2 = 'Ptprz1'
$ws.Range('D11').Value2 = 'ECs'
$ws.Range('E11').Value2 = [double]"3"
$ws.Range('F11').Value2 = [double]"1"
$ws.Range('G11').Value2 = [double]"15.96019966666667"
$ws.Range('H11').Value2 = [double]"47.880599"
$ws.Range('I11').Value2 = [double]"0.6184839868322428"
$ws.Range('J11').Value2 = [double]"0.6184839868322429"
$ws.Range('K11').Value2 = [double]"1"
$ws.Range('L11').Value2 = [double]"0.3333333333333333"
$ws.Range('M11').Value2 = [double]"0.04936366666666667"
$ws.Range('N11').Value2 = [double]"0.148091"
$ws.Range('O11').Value2 = [double]"0.04616170608573571"
$ws.Range('P11').Value2 = [double]"0.0461617060857357"
$ws.Range('Q11').Value2 = [double]"0.7878539762787778"
$ws.Range('R11').Value2 = [double]"7.090685786509001"
$ws.Range('S11').Value2 = [double]"0.02855027601888403"
$ws.Range('T11').Value2 = [double]"0.02855027601888403"

# Row 12
$ws.Range('A12').Value2 = 'MuSCs'
$ws.Range('B12').Value2 = 'Ncam1'
$ws.Range('C12').Value2 = 'Ptprz1'
$ws.Range('D12').Value2 = 'FAPs'
$ws.Range('E12').Value2 = [double]"3"
$ws.Range('F12').Value2 = [double]"1"
$ws.Range('G12').Value2 = [double]"15.96019966666667"
$ws.Range('H12').Value2 = [double]"47.880599"
$ws.Range('I12').Value2 = [double]"0.6184839868322428"
$ws.Range('J12').Value2 = [double]"0.6184839868322429"
$ws.Range('K12').Value2 = [double]"1"
$ws.Range('L12').Value2 = [double]"0.3333333333333333"
$ws.Range('M12').Value2 = [double]"0.008616666666666667"
$ws.Range('N12').Value2 = [double]"0.02585"
$ws.Range('O12').Value2 = [double]"0.008057748967298944"
$ws.Range('P12').Value2 = [double]"0.008057748967298944"
$ws.Range('Q12').Value2 = [double]"0.1375237204611111"
$ws.Range('R12').Value2 = [double]"1.23771348415"
$ws.Range('S12').Value2 = [double]"0.004983588706188438"
$ws.Range('T12').Value2 = [double]"0.004983588706188439"

# Row 13
$ws.Range('A13').Value2 = 'MuSCs'
$ws.Range('B13').Value2 = 'Ncam1'
$ws.Range('C13').Value2 = 'Ptprz1'
$ws.Range('D13').Value2 = 'MuSCs'
$ws.Range('E13').Value2 = [double]"3"
$ws.Range('F13').Value2 = [double]"1"
$ws.Range('G13').Value2 = [double]"15.96019966666667"
$ws.Range('H13').Value2 = [double]"47.880599"
$ws.Range('I13').Value2 = [double]"0.6184839868322428"
$ws.Range('J13').Value2 = [double]"0.6184839868322429"
$ws.Range('K13').Value2 = [double]"3"
$ws.Range('L13').Value2 = [double]"1"
$ws.Range('M13').Value2 = [double]"1.011383666666666"
$ws.Range('N13').Value2 = [double]"3.034151"
$ws.Range('O13').Value2 = [double]"0.9457805449469654"
$ws.Range('P13').Value2 = [double]"0.9457805449469653"
$ws.Range('Q13').Value2 = [double]"16.14188525960544"
$ws.Range('R13').Value2 = [double]"145.276967336449"
$ws.Range('S13').Value2 = [double]"0.5849501221071703"
$ws.Range('T13').Value2 = [double]"0.5849501221071705"

# Row 14
$ws.Range('A14').Value2 = 'Neutrophils'
$ws.Range('B14').Value2 = 'Ncam1'
$ws.Range('C14').Value2 = 'Ptprz1'
$ws.Range('D14').Value2 = 'ECs'
$ws.Range('E14').Value2 = [double]"3"
$ws.Range('F14').Value2 = [double]"1"
$ws.Range('G14').Value2 = [double]"0.3268106666666666"
$ws.Range('H14').Value2 = [double]"0.980432"
$ws.Range('I14').Value2 = [double]"0.01266445083901121"
$ws.Range('J14').Value2 = [double]"0.01266445083901121"
$ws.Range('K14').Value2 = [double]"1"
$ws.Range('L14').Value2 = [double]"0.3333333333333333"
$ws.Range('M14').Value2 = [double]"0.04936366666666667"
$ws.Range('N14').Value2 = [double]"0.148091"
$ws.Range('O14').Value2 = [double]"0.04616170608573571"
$ws.Range('P14').Value2 = [double]"0.0461617060857357"
$ws.Range('Q14').Value2 = [double]"0.01613257281244444"
$ws.Range('R14').Value2 = [double]"0.145193155312"
$ws.Range('S14').Value2 = [double]"0.0005846126573676846"
$ws.Range('T14').Value2 = [double]"0.0005846126573676846"

# Row 15
$ws.Range('A15').Value2 = 'Neutrophils'
$ws.Range('B15').Value2 = 'Ncam1'
$ws.Range('C15').Value2 = 'Ptprz1'
$ws.Range('D15').Value2 = 'FAPs'
$ws.Range('E15').Value2 = [double]"3"
$ws.Range('F15').Value2 = [double]"1"
$ws.Range('G15').Value2 = [double]"0.3268106666666666"
$ws.Range('H15').Value2 = [double]"0.980432"
$ws.Range('I15').Value2 = [double]"0.01266445083901121"
$ws.Range('J15').Value2 = [double]"0.01266445083901121"
$ws.Range('K15').Value2 = [double]"1"
$ws.Range('L15').Value2 = [double]"0.3333333333333333"
$ws.Range('M15').Value2 = [double]"0.008616666666666667"
$ws.Range('N15').Value2 = [double]"0.02585"
$ws.Range('O15').Value2 = [double]"0.008057748967298944"
$ws.Range('P15').Value2 = [double]"0.008057748967298944"
$ws.Range('Q15').Value2 = [double]"0.002816018577777777"
$ws.Range('R15').Value2 = [double]"0.0253441672"
$ws.Range('S15').Value2 = [double]"0.0001020469656694508"
$ws.Range('T15').Value2 = [double]"0.0001020469656694508"

# Row 16
$ws.Range('A16').Value2 = 'Neutrophils'
$ws.Range('B16').Value2 = 'Ncam1'
$ws.Range('C16').Value2 = 'Ptprz1'
$ws.Range('D16').Value2 = 'MuSCs'
$ws.Range('E16').Value2 = [double]"3"
$ws.Range('F16').Value2 = [double]"1"
$ws.Range('G16').Value2 = [double]"0.3268106666666666"
$ws.Range('H16').Value2 = [double]"0.980432"
$ws.Range('I16').Value2 = [double]"0.01266445083901121"
$ws.Range('J16').Value2 = [double]"0.01266445083901121"
$ws.Range('K16').Value2 = [double]"3"
$ws.Range('L16').Value2 = [double]"1"
$ws.Range('M16').Value2 = [double]"1.011383666666666"
$ws.Range('N16').Value2 = [double]"3.034151"
$ws.Range('O16').Value2 = [double]"0.9457805449469654"
$ws.Range('P16').Value2 = [double]"0.9457805449469653"
$ws.Range('Q16').Value2 = [double]"0.330530970359111"
$ws.Range('R16').Value2 = [double]"2.974778733232"
$ws.Range('S16').Value2 = [double]"0.01197779121597408"
$ws.Range('T16').Value2 = [double]"0.01197779121597408"

# Row 17
$ws.Range('A17').Value2 = 'Resolving-Mac'
$ws.Range('B17').Value2 = 'Ncam1'
$ws.Range('C17').Value2 = 'Ptprz1'
$ws.Range('D17').Value2 = 'ECs'
$ws.Range('E17').Value2 = [double]"1"
$ws.Range('F17').Value2 = [double]"0.3333333333333333"
$ws.Range('G17').Value2 = [double]"0.008676333333333333"
$ws.Range('H17').Value2 = [double]"0.026029"
$ws.Range('I17').Value2 = [double]"0.0003362221866367304"
$ws.Range('J17').Value2 = [double]"0.0003362221866367304"
$ws.Range('K17').Value2 = [double]"1"
$ws.Range('L17').Value2 = [double]"0.3333333333333333"
$ws.Range('M17').Value2 = [double]"0.04936366666666667"
$ws.Range('N17').Value2 = [double]"0.148091"
$ws.Range('O17').Value2 = [double]"0.04616170608573571"
$ws.Range('P17').Value2 = [double]"0.0461617060857357"
$ws.Range('Q17').Value2 = [double]"0.0004282956265555555"
$ws.Range('R17').Value2 = [double]"0.003854660639"
$ws.Range('S17').Value2 = [double]"1.552058975902812E-05"
$ws.Range('T17').Value2 = [double]"1.552058975902812E-05"

# Row 18
$ws.Range('A18').Value2 = 'Resolving-Mac'
$ws.Range('B18').Value2 = 'Ncam1'
$ws.Range('C18').Value2 = 'Ptprz1'
$ws.Range('D18').Value2 = 'FAPs'
$ws.Range('E18').Value2 = [double]"1"
$ws.Range('F18').Value2 = [double]"0.3333333333333333"
$ws.Range('G18').Value2 = [double]"0.008676333333333333"
$ws.Range('H18').Value2 = [double]"0.026029"
$ws.Range('I18').Value2 = [double]"0.0003362221866367304"
$ws.Range('J18').Value2 = [double]"0.0003362221866367304"
$ws.Range('K18').Value2 = [double]"1"
$ws.Range('L18').Value2 = [double]"0.3333333333333333"
$ws.Range('M18').Value2 = [double]"0.008616666666666667"
$ws.Range('N18').Value2 = [double]"0.02585"
$ws.Range('O18').Value2 = [double]"0.008057748967298944"
$ws.Range('P18').Value2 = [double]"0.008057748967298944"
$ws.Range('Q18').Value2 = [double]"7.476107222222222E-05"
$ws.Range('R18').Value2 = [double]"0.00067284965"
$ws.Range('S18').Value2 = [double]"2.709193977155107E-06"
$ws.Range('T18').Value2 = [double]"2.709193977155108E-06"

# Row 19
$ws.Range('A19').Value2 = 'Resolving-Mac'
$ws.Range('B19').Value2 = 'Ncam1'
$ws.Range('C19').Value2 = 'Ptprz1'
$ws.Range('D19').Value2 = 'MuSCs'
$ws.Range('E19').Value2 = [double]"1"
$ws.Range('F19').Value2 = [double]"0.3333333333333333"
$ws.Range('G19').Value2 = [double]"0.008676333333333333"
$ws.Range('H19').Value2 = [double]"0.026029"
$ws.Range('I19').Value2 = [double]"0.0003362221866367304"
$ws.Range('J19').Value2 = [double]"0.0003362221866367304"
$ws.Range('K19').Value2 = [double]"3"
$ws.Range('L19').Value2 = [double]"1"
$ws.Range('M19').Value2 = [double]"1.011383666666666"
$ws.Range('N19').Value2 = [double]"3.034151"
$ws.Range('O19').Value2 = [double]"0.9457805449469654"
$ws.Range('P19').Value2 = [double]"0.9457805449469653"
$ws.Range('Q19').Value2 = [double]"0.008775101819888886"
$ws.Range('R19').Value2 = [double]"0.07897591637899999"
$ws.Range('S19').Value2 = [double]"0.0003179924029005472"
$ws.Range('T19').Value2 = [double]"0.0003179924029005472"
